$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as TEXT (preserving formats like "1.00" or "0.999"),
# then strip the quote-prefix style so the cell keeps the default (unstyled) look.
function Set-TextCell($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" '58.948.93'
Set-TextCell "E2" '  -0.38%  '
Set-TextCell "D3" '2.509.02'
Set-TextCell "E3" '  -0.77%  '
Set-TextCell "E4" '  +0.22%  '
Set-TextCell "D5" '533.46'
Set-TextCell "E5" '  -0.46%  '
Set-TextCell "D6" '135.41'
Set-TextCell "E6" '  -2.05%  '
Set-TextCell "D7" '0.999'
Set-TextCell "E7" '  +0.19%  '
Set-TextCell "E8" '  -0.18%  '
Set-TextCell "E9" '  -0.79%  '
Set-TextCell "E10" '  -1.42%  '
Set-TextCell "E12" '  -0.51%  '
Set-TextCell "D13" '2.953.93'
Set-TextCell "E13" '  -0.19%  '
Set-TextCell "D14" '58.838.06'
Set-TextCell "E14" '  -0.42%  '
Set-TextCell "D15" '22.81'
Set-TextCell "E15" '  -2.10%  '
Set-TextCell "E16" '  -1.23%  '
Set-TextCell "D17" '2.500.17'
Set-TextCell "E17" '  -0.30%  '
Set-TextCell "D18" '11.03'
Set-TextCell "E18" '  -0.78%  '
Set-TextCell "D19" '4.26'
Set-TextCell "E19" '  -0.48%  '
Set-TextCell "D20" '323.02'
Set-TextCell "E20" '  -0.77%  '
Set-TextCell "E21" '  -0.45%  '
Set-TextCell "D23" '65.06'
Set-TextCell "E23" '  +0.04%  '
Set-TextCell "E24" '  -0.20%  '
Set-TextCell "E25" '  -1.68%  '
Set-TextCell "E26" '  -1.31%  '
Set-TextCell "E27" '  -0.83%  '
Set-TextCell "D28" '0.0₃0761'
Set-TextCell "E28" '  -2.19%  '
Set-TextCell "D29" '6.49'
Set-TextCell "E29" '  -3.57%  '
Set-TextCell "D31" '168.70'
Set-TextCell "E31" '  +0.29%  '
Set-TextCell "E32" '  +0.11%  '
Set-TextCell "E33" '  -5.24%  '
Set-TextCell "B34" 'EthereumClassic'
Set-TextCell "C34" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell "D34" '18.43'
Set-TextCell "E34" '  -0.96%  '
Set-TextCell "B35" 'ImmutableX'
Set-TextCell "C35" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D35" '1.36'
Set-TextCell "E35" '  -4.23%  '
Set-TextCell "D36" '4.04'
Set-TextCell "E36" '  -2.00%  '
Set-TextCell "E37" '  -2.85%  '
Set-TextCell "D38" '3.57'
Set-TextCell "E38" '  -2.26%  '
Set-TextCell "D39" '0.798'
Set-TextCell "E39" '  -4.26%  '
Set-TextCell "D40" '281.67'
Set-TextCell "E40" '  -0.33%  '
Set-TextCell "E41" '  -0.11%  '
Set-TextCell "E42" '  -0.36%  '
Set-TextCell "D43" '5.00'
Set-TextCell "E43" '  -4.92%  '
Set-TextCell "E44" '  -0.87%  '
Set-TextCell "E45" '  +0.47%  '
Set-TextCell "E46" '  -0.74%  '
Set-TextCell "E47" '  -2.70%  '
Set-TextCell "E48" '  -2.76%  '
Set-TextCell "D49" '17.26'
Set-TextCell "E49" '  -0.93%  '
Set-TextCell "D50" '1.757.01'
Set-TextCell "E50" '  -0.37%  '
Set-TextCell "D51" '0.983'
Set-TextCell "E51" '  -0.59%  '
